# Auto update nse_indices_1 outputs
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: Average Market Change (%) text value
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$cell = $wsOverview.Cells.Item(3, 2)
$cell.NumberFormat = "@"
$cell.Value = "4.89%"

# ---------------------------------------------------------------------
# Summary sheet: MTD % table re-ranked (row 3..7 index names + values)
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Cells.Item(3, 1).Value = "Nifty Midcap 100"
$wsSummary.Cells.Item(3, 2).Value = 6.13
$wsSummary.Cells.Item(4, 1).Value = "Nifty 50"
$wsSummary.Cells.Item(4, 2).Value = 5.45
$wsSummary.Cells.Item(5, 1).Value = "Nifty 200"
$wsSummary.Cells.Item(5, 2).Value = 5.2
$wsSummary.Cells.Item(6, 1).Value = "Nifty 100"
$wsSummary.Cells.Item(6, 2).Value = 4.98
$wsSummary.Cells.Item(7, 1).Value = "Nifty Midcap 150"
$wsSummary.Cells.Item(7, 2).Value = 4.9

# ---------------------------------------------------------------------
# Index Close sheet
# ---------------------------------------------------------------------
$wsClose = $wb.Worksheets.Item("Index Close")
$wsClose.Cells.Item(2, 7).Value = 60310.1484375

$wsClose.Cells.Item(7, 2).Value = 25898.55078125
$wsClose.Cells.Item(7, 3).Value = 67939.5
$wsClose.Cells.Item(7, 4).Value = 26402.44921875
$wsClose.Cells.Item(7, 6).Value = 17089.75
$wsClose.Cells.Item(7, 7).Value = 59578.05078125
$wsClose.Cells.Item(7, 9).Value = 23550.849609375
$wsClose.Cells.Item(7, 11).Value = 23550.849609375

# ---------------------------------------------------------------------
# MTD % sheet
# ---------------------------------------------------------------------
$wsMtd = $wb.Worksheets.Item("MTD %")
$wsMtd.Cells.Item(2, 7).Value = 6.13

$wsMtd.Cells.Item(7, 2).Value = 5.17
$wsMtd.Cells.Item(7, 3).Value = 1.9
$wsMtd.Cells.Item(7, 4).Value = 4.64
$wsMtd.Cells.Item(7, 6).Value = 6.59
$wsMtd.Cells.Item(7, 7).Value = 4.84
$wsMtd.Cells.Item(7, 9).Value = 3.69
$wsMtd.Cells.Item(7, 11).Value = 3.69

# ---------------------------------------------------------------------
# DoD% sheet
# ---------------------------------------------------------------------
$wsDod = $wb.Worksheets.Item("DoD%")
$wsDod.Cells.Item(2, 7).Value = 1.2

$wsDod.Cells.Item(6, 2).Value = 0.5
$wsDod.Cells.Item(6, 3).Value = 1.01
$wsDod.Cells.Item(6, 4).Value = 0.58
$wsDod.Cells.Item(6, 6).Value = 1.01
$wsDod.Cells.Item(6, 7).Value = 1.07
$wsDod.Cells.Item(6, 9).Value = 0.74
$wsDod.Cells.Item(6, 11).Value = 0.74

$wsDod.Cells.Item(7, 2).Value = 0
$wsDod.Cells.Item(7, 3).Value = 0
$wsDod.Cells.Item(7, 4).Value = 0
$wsDod.Cells.Item(7, 6).Value = 0
$wsDod.Cells.Item(7, 7).Value = 0
$wsDod.Cells.Item(7, 9).Value = 0
$wsDod.Cells.Item(7, 11).Value = 0

# ---------------------------------------------------------------------
# Daily Movers sheet: Top 3 Gainers / Top 3 Losers text lists
# ---------------------------------------------------------------------
$wsMovers = $wb.Worksheets.Item("Daily Movers")
$wsMovers.Cells.Item(2, 2).Value = "Nifty Midcap 100, Nifty Midcap 150, Nifty Midcap 50"
$wsMovers.Cells.Item(2, 3).Value = "Nifty 50, Nifty 100, Nifty 200"

$wsMovers.Cells.Item(6, 2).Value = "Nifty Midcap 100, Nifty Next 50, Nifty Midcap 50"
$wsMovers.Cells.Item(6, 3).Value = "Nifty Midcap 150, Nifty 200, Nifty500 Multicap 50:25:25"

$wsMovers.Cells.Item(7, 2).Value = "Nifty Midcap 150, Nifty500 Multicap 50:25:25, Nifty 200"
$wsMovers.Cells.Item(7, 3).Value = "Nifty 50, Nifty Next 50, Nifty 100"
